$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the "Final" (Post Treatment outcome) column values in D2:D7.
# D2 and D3 mirror column C (A little stressful), while D4:D7 are
# "Moderately stressful" to reflect the final parent outcome measurements.
$ws.Range("D2").Value = "A little stressful"
$ws.Range("D3").Value = "A little stressful"
$ws.Range("D4").Value = "Moderately stressful"
$ws.Range("D5").Value = "Moderately stressful"
$ws.Range("D6").Value = "Moderately stressful"
$ws.Range("D7").Value = "Moderately stressful"

# Column D needs to widen to fit the new content
$ws.Columns.Item(4).ColumnWidth = 19

# Update the active selection to D8 (just below the new data)
$ws.Range("D8").Select()
